# Refresh the crypto price/volume snapshot (GitHub Actions data pull).
# Column layout: A=index, B=Coin, C=Link, D=Price, E=Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "69.174.69"
$ws.Cells.Item(2, 5).Value = "  +1.31%  "
# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.932.84"
$ws.Cells.Item(3, 5).Value = "  +0.20%  "
# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "493.68"
$ws.Cells.Item(5, 5).Value = "  +1.31%  "
# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "147.37"
$ws.Cells.Item(6, 5).Value = "  -0.97%  "
# Row 7
$ws.Cells.Item(7, 5).Value = "  -1.22%  "
# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.04%  "
# Row 9
$ws.Cells.Item(9, 5).Value = "  -0.42%  "
# Row 10
$ws.Cells.Item(10, 5).Value = "  +4.17%  "
# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.97%  "
# Row 12
$ws.Cells.Item(12, 5).Value = "  +0.51%  "
# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "10.45"
$ws.Cells.Item(13, 5).Value = "  -2.40%  "
# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.561.33"
$ws.Cells.Item(14, 5).Value = "  +0.26%  "
# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.935.64"
$ws.Cells.Item(15, 5).Value = "  +0.16%  "
# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "14.28"
$ws.Cells.Item(16, 5).Value = "  -3.36%  "
# Row 17
$ws.Cells.Item(17, 5).Value = "  -0.83%  "
# Row 18
$ws.Cells.Item(18, 5).Value = "  +4.42%  "
# Row 19
$ws.Cells.Item(19, 5).Value = "  -1.02%  "
# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "69.227.52"
$ws.Cells.Item(20, 5).Value = "  +1.29%  "
# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "438.91"
$ws.Cells.Item(21, 5).Value = "  -0.65%  "
# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "3.45"
$ws.Cells.Item(22, 5).Value = "  +0.25%  "
# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "14.52"
$ws.Cells.Item(23, 5).Value = "  -2.92%  "
# Row 24
$ws.Cells.Item(24, 2).Value = "Litecoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "88.71"
$ws.Cells.Item(24, 5).Value = "  +0.09%  "
# Row 25
$ws.Cells.Item(25, 2).Value = "RenderToken"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "12.24"
$ws.Cells.Item(25, 5).Value = "  +9.76%  "
# Row 26
$ws.Cells.Item(26, 5).Value = "  +4.93%  "
# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "11.13"
$ws.Cells.Item(27, 5).Value = "  -2.83%  "
# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "37.12"
$ws.Cells.Item(28, 5).Value = "  -4.18%  "
# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.65"
$ws.Cells.Item(29, 5).Value = "  -3.87%  "
# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "701.56"
$ws.Cells.Item(30, 5).Value = "  -2.34%  "
# Row 31
$ws.Cells.Item(31, 2).Value = "Cosmos"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "13.39"
$ws.Cells.Item(31, 5).Value = "  -1.53%  "
# Row 32
$ws.Cells.Item(32, 2).Value = "Hedera"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.131"
$ws.Cells.Item(32, 5).Value = "  +0.06%  "
# Row 33
$ws.Cells.Item(33, 5).Value = "  -0.02%  "
# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.470"
$ws.Cells.Item(34, 5).Value = "  +18.81%  "
# Row 35
$ws.Cells.Item(35, 5).Value = "  -1.37%  "
# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "6.13"
$ws.Cells.Item(36, 5).Value = "  -0.55%  "
# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "61.92"
$ws.Cells.Item(37, 5).Value = "  +3.60%  "
# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "40.72"
$ws.Cells.Item(38, 5).Value = "  -3.40%  "
# Row 39
$ws.Cells.Item(39, 5).Value = "  +0.56%  "
# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.999"
$ws.Cells.Item(40, 5).Value = "  -0.13%  "
# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.00%  "
# Row 42
$ws.Cells.Item(42, 5).Value = "  +1.27%  "
# Row 43
$ws.Cells.Item(43, 5).Value = "  -2.02%  "
# Row 44
$ws.Cells.Item(44, 5).Value = "  -3.73%  "
# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.99"
$ws.Cells.Item(45, 5).Value = "  +1.48%  "
# Row 46
$ws.Cells.Item(46, 5).Value = "  +0.43%  "
# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "3.39"
$ws.Cells.Item(47, 5).Value = "  +7.58%  "
# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0₆0355"
$ws.Cells.Item(48, 5).Value = "  -0.81%  "
# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "3.01"
$ws.Cells.Item(49, 5).Value = "  +6.27%  "
# Row 50
$ws.Cells.Item(50, 5).Value = "  -1.26%  "
# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "144.44"
